$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "corrected error in fixed recourse data": update solve time (D), gap (C), and objective (B) values
$ws.Range("D2").Value = 3671.16128473

$ws.Range("C3").Value = 0.3524663646928933
$ws.Range("D3").Value = 3623.585527484

$ws.Range("D4").Value = 3675.185769438

$ws.Range("C5").Value = 0.42494411658262016
$ws.Range("D5").Value = 3620.7815772

$ws.Range("D6").Value = 3625.439816565

$ws.Range("D7").Value = 3621.975197309

$ws.Range("C8").Value = 0.46065307544108214
$ws.Range("D8").Value = 3628.08236608

$ws.Range("B9").Value = -296.7607524218697
$ws.Range("C9").Value = 0.4568009651124302
$ws.Range("D9").Value = 3622.792508465

$ws.Range("D10").Value = 3624.029671429

$ws.Range("D11").Value = 3716.806027163
